$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Kartik Tyagi"
$ws.Name = "Kartik Tyagi"

# Insert a new column before column A ("matchNo") - everything shifts right
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - update / fill in values (column B onward already shifted from insert)
# Numeric-looking values are prefixed with a leading apostrophe so they are
# stored as literal text (matching the source "numberStoredAsText" data),
# instead of being auto-converted to real numbers by Excel.
$ws.Range("A2").Value = "43rd"
$ws.Range("B2").Value = "Rajasthan Royals"
$ws.Range("C2").Value = "Kartik Tyagi"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'1"
$ws.Range("F2").Value = "'1"
$ws.Range("G2").Value = "'0"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'100.00"
$ws.Range("J2").Value = "Royal Challengers Bangalore"
$ws.Range("K2").Value = "Dubai (DSC)"
$ws.Range("L2").Value = "September 29"
$ws.Range("M2").Value = "RCB won by 7 wickets (with 17 balls remaining)"

# Row 3 - new row of data
$ws.Range("A3").Value = "32nd"
$ws.Range("B3").Value = "Rajasthan Royals"
$ws.Range("C3").Value = "Kartik Tyagi"
$ws.Range("D3").Value = "b Arshdeep Singh"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'3"
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'33.33"
$ws.Range("J3").Value = "Punjab Kings"
$ws.Range("K3").Value = "Dubai (DSC)"
$ws.Range("L3").Value = "September 21"
$ws.Range("M3").Value = "Royals won by 2 runs"
